$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "3.3V Supply" extra row: Current Sensor x 6 (columns D:E) ---
$ws.Range("D2").Value = "Current Sensor x 6"
$ws.Range("E2").Formula = "=25*6"

# --- New "5V Supply" block (columns G:H) ---
$ws.Range("G1").Value = "5V Supply"
$ws.Range("H1").Value = "Max Current Draw (mA)"
$ws.Range("G1:H1").Font.Bold = $true

$ws.Range("G2").Value = "CAN Transceiver"

# --- New "3.3V Supply" extra row: External ADC (columns D:E) ---
$ws.Range("D3").Value = "External ADC"
$ws.Range("E3").Value = 2.5

# --- New rows appended to the 12V Supply block (columns A:B) ---
$ws.Range("A5").Value = "Gigavac contactor"
$ws.Range("B5").Value = 170

$ws.Range("A6").Value = "Omron Relay x 2"
$ws.Range("B6").Formula = "=100*2"

# Column widths widened by Excel's "best fit" as a side-effect of the new,
# wider text in these columns (closest achievable values in this engine).
$ws.Columns.Item(4).ColumnWidth = 14.3671875
$ws.Columns.Item(7).ColumnWidth = 12.78125
$ws.Columns.Item(8).ColumnWidth = 19.3515625

# Update selection to match the authored state
$ws.Range("F17").Select()
